$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 64331
$ws.Range("B2").Value = "Agatha Montenegro"
$ws.Range("C2").Value = "Vendas"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 45103
$ws.Range("G2").Value = 6158.58

# Row 3
$ws.Range("A3").Value = 75313
$ws.Range("B3").Value = "Liz Rios"
$ws.Range("C3").Value = "Atendimento ao Cliente"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 45081
$ws.Range("G3").Value = 4870.06

# Row 4
$ws.Range("A4").Value = 57518
$ws.Range("B4").Value = "Eduarda Aparecida"
$ws.Range("C4").Value = "Engenharia"
$ws.Range("D4").Value = "Doenca"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 45099
$ws.Range("G4").Value = 8643.360000000001

# Row 5
$ws.Range("A5").Value = 76551
$ws.Range("B5").Value = "Sr. Léo Barros"
$ws.Range("C5").Value = "Atendimento ao Cliente"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 45095
$ws.Range("G5").Value = 5341.8

# Row 6
$ws.Range("A6").Value = 87866
$ws.Range("B6").Value = "Otto Farias"
$ws.Range("C6").Value = "P&D"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 45103
$ws.Range("G6").Value = 8712.200000000001

# Row 7
$ws.Range("A7").Value = 5248
$ws.Range("B7").Value = "Aylla Fernandes"
$ws.Range("C7").Value = "Marketing"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 45085
$ws.Range("G7").Value = 2618.59

# Row 8
$ws.Range("A8").Value = 38176
$ws.Range("B8").Value = "Maria Cecília Abreu"
$ws.Range("C8").Value = "TI"
$ws.Range("D8").Value = "Consulta medica"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 45094
$ws.Range("G8").Value = 9665.309999999999

# Row 9
$ws.Range("A9").Value = 90512
$ws.Range("B9").Value = "Maria Vitória Pires"
$ws.Range("C9").Value = "Operacoes"
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 45082
$ws.Range("G9").Value = 5593.33

# Row 10
$ws.Range("A10").Value = 34694
$ws.Range("B10").Value = "Murilo Souza"
$ws.Range("C10").Value = "Engenharia"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 45090
$ws.Range("G10").Value = 6271.71

# Row 11
$ws.Range("A11").Value = 87907
$ws.Range("B11").Value = "Luara Correia"
$ws.Range("C11").Value = "Atendimento ao Cliente"
$ws.Range("D11").Value = "Doenca"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 45103
$ws.Range("G11").Value = 8015.5

$wb.Save()
